# Fruta / hortaliza, semanal
# Insert two new weekly price rows right after the existing row 3
# (pushing the former rows 4-28 down to rows 6-30), then populate the
# two new rows (4 and 5) with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 4 - shifts old rows 4:28 to 6:30
$ws.Rows("4:5").Insert()

# ---- New row 4: Packham's Triumph ----
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44715
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104005
$ws.Range("J4").Value = "Pera"
$ws.Range("K4").Value = "Packham's Triumph"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18

# ---- New row 5: Winter Nelis ----
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44715
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104005
$ws.Range("J5").Value = "Pera"
$ws.Range("K5").Value = "Winter Nelis"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 972
$ws.Range("T5").Value = 18
